# Applies the "completed functional simulation screenshots and working on
# some documentations" update to the Activity Log workbook.
#
# Order of operations matters: the shared-strings table is built up in the
# exact order new text values are first entered (mirrors how Excel appends
# to sst on save), so cells are written in that sequence.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Activity Log - Part 1")
$ws2 = $wb.Worksheets.Item("Activity Log - Part 2")
$ws3 = $wb.Worksheets.Item("Activity Log - Part 3")

# ---------------------------------------------------------------------
# Sheet 1 - "Activity Log - Part 1"
# ---------------------------------------------------------------------

$rows1 = @(
    @{ Row=6;  Digits=4794; Date=43919; Start=0.88888888888888884; End=0.90972222222222221; Desc="Started reading about the implementation of Arithmetic Unit" },
    @{ Row=7;  Digits=4794; Date=43919; Start=0.91319444444444453; End=0.94791666666666663; Desc="Setting up work environment, git and modelsim" },
    @{ Row=8;  Digits=4794; Date=43920; Start=0.86805555555555547; End=0.90625;              Desc="Worked on full adder implementation" },
    @{ Row=9;  Digits=4794; Date=43920; Start=0.90625;              End=0.95138888888888884; Desc="Worked on ripple adder implementation" },
    @{ Row=10; Digits=4794; Date=43921; Start=0.86111111111111116; End=0.93055555555555547; Desc="Worked on arithmetic unit implementing Adder, Zero, ExtWord MUX, AltB  and AltBu" },
    @{ Row=11; Digits=4794; Date=43922; Start=0.60416666666666663; End=0.66666666666666663; Desc="Help debug arithmetic unit (errors with sign extension) Output ExtWord was not matching with test bench values" },
    @{ Row=12; Digits=4794; Date=43923; Start=0.72916666666666663; End=0.77083333333333337; Desc="Screenshots of waves of functional simulation" },
    @{ Row=13; Digits=4794; Date=43923; Start=0.81944444444444453; End=0.86458333333333337; Desc="Helping with screenshots of timing simulations" },
    @{ Row=14; Digits=4794; Date=43923; Start=0.86458333333333337; End=0.94791666666666663; Desc="Start working on documentations, screenshot descriptions, etc." },
    @{ Row=15; Digits=4794; Date=43924; Start=0.71875;              End=0.76388888888888884; Desc="Writing up report and proofreading" },
    @{ Row=16; Digits=4794; Date=43926; Start=0.61111111111111105; End=0.67708333333333337; Desc="Cleaning up documentations and finishing up" },
    @{ Row=17; Digits=4794; Date=43926; Start=0.67708333333333337; End=0.72916666666666663; Desc="Adding anotations to pdf" },
    @{ Row=18; Digits=4794; Date=43926; Start=0.88541666666666663; End=0.93055555555555547; Desc="Added table of contents, formatting and submitting" }
)

foreach ($r in $rows1) {
    $n = $r.Row
    $ws1.Range("B$n").Value = $r.Digits
    $ws1.Range("C$n").Value = $r.Date
    $ws1.Range("D$n").Value = $r.Start
    $ws1.Range("E$n").Value = $r.End
    $ws1.Range("G$n").Value = $r.Desc
}

# Header block: group number + student name
$ws1.Range("B3").Value = "G54"
$ws1.Range("B1").Value = "Yoel Yonata"
$ws1.Range("B2").Value = 301304794

# ---------------------------------------------------------------------
# Sheet 2 - "Activity Log - Part 2"
# ---------------------------------------------------------------------

$ws2.Range("G7").Value = "Worked on SLL, SRL and SRA"
$ws2.Range("G8").Value = "Worked on Execution Unit"
$ws2.Range("G6").Value = "Started reading on project part 2 and how to implement"
$ws2.Range("G9").Value = "Debugging SRA, sign extension was not working"

$ws2.Range("B6").Value = 4794
$ws2.Range("C6").Value = 43931
$ws2.Range("D6").Value = 0.88541666666666663
$ws2.Range("E6").Value = 0.92708333333333337

$ws2.Range("B7").Value = 4794
$ws2.Range("C7").Value = 43931
$ws2.Range("D7").Value = 0.92708333333333337
$ws2.Range("E7").Value = 0.98958333333333337

$ws2.Range("B8").Value = 4794
$ws2.Range("C8").Value = 43932
$ws2.Range("D8").Value = 0.86458333333333337
$ws2.Range("E8").Value = 0.90972222222222221

$ws2.Range("B9").Value = 4794
$ws2.Range("C9").Value = 43932
$ws2.Range("D9").Value = 0.90972222222222221
$ws2.Range("E9").Value = 0.99305555555555547

$ws2.Range("B3").Value = "G54"
$ws2.Range("B1").Value = "Yoel Yonata"
$ws2.Range("B2").Value = 301304794

# ---------------------------------------------------------------------
# View state: sheet 2 becomes the active tab, both sheet 1 & 2 zoomed to
# 70%, and selections move to reflect where the student left off.
# ---------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("G27").Select()
$excel.ActiveWindow.Zoom = 70

$ws2.Activate()
$ws2.Range("G10").Select()
$excel.ActiveWindow.Zoom = 70

$ws3.Activate()
$ws3.Range("G12").Select()
